$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All data cells in this sheet are stored as text (inline strings) in the
# original workbook, including values that look numeric (e.g. "0.4609",
# "1.000", "26.884.93"). Force each target cell to Text format before
# writing so Excel does not reinterpret the string as a number/date and
# round, truncate, or reformat it.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.884.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.805.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.88"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4609"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3732"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07369"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8738"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.60%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.752.16"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.95%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.351"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.516"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07047"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "90.75"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.33%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008734"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.898.11"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.313"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.906.25"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -7.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.16"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.24%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.142"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -9.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.294"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.78"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08897"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7676"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.159"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.470"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.114"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05248"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.403"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +4.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5336"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.217"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.895"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.34%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1659"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.558"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5050"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.30"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.62"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.26%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.05%  "
